$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Import "lop" (class) data from Excel -------------------------------
# Row 2 - note: D (makhoa) is entered before C (nienkhoa) to reproduce the
# original authoring order (affects shared-string table order).
$ws.Cells.Item(2, 1).Value = "DHCN3A"
$ws.Cells.Item(2, 2).Value = "Đại học công nghệ 3A"
$ws.Cells.Item(2, 4).Value = "CNTT-TCKGM"
$ws.Cells.Item(2, 3).Value = "2016-2020"

# Row 3
$ws.Cells.Item(3, 1).Value = "DHCN3B"
$ws.Cells.Item(3, 2).Value = "Đại học công nghệ 3B"
$ws.Cells.Item(3, 3).Value = "2016-2020"
$ws.Cells.Item(3, 4).Value = "CNTT-TCKGM"

# Row 4
$ws.Cells.Item(4, 1).Value = "DHCN4"
$ws.Cells.Item(4, 2).Value = "Đại học công nghệ 4"
$ws.Cells.Item(4, 3).Value = "2017-2021"
$ws.Cells.Item(4, 4).Value = "CNTT-TCKGM"

# Row 5
$ws.Cells.Item(5, 1).Value = "DHCN5"
$ws.Cells.Item(5, 2).Value = "Đại học công nghệ 5"
$ws.Cells.Item(5, 3).Value = "2018-2022"
$ws.Cells.Item(5, 4).Value = "CNTT-TCKGM"

# Row 6
$ws.Cells.Item(6, 1).Value = "DTVT3"
$ws.Cells.Item(6, 2).Value = "Điện tử viễn thông 3"
$ws.Cells.Item(6, 4).Value = "DTVT"
$ws.Cells.Item(6, 3).Value = "2016-2020"

# Row 7
$ws.Cells.Item(7, 1).Value = "DTVT4"
$ws.Cells.Item(7, 2).Value = "Điện tử viễn thông 4"
$ws.Cells.Item(7, 3).Value = "2016-2020"
$ws.Cells.Item(7, 4).Value = "DTVT"

# --- Header formatting: center both horizontally and vertically --------
# Build the combined (horizontal+vertical center) format on a scratch
# single cell first, then copy/paste the format onto the header row. A
# single property assignment on a multi-cell range leaves a stray
# intermediate style behind, so we stage it on one cell instead.
$tmpl = $ws.Range("Z1000")
$tmpl.HorizontalAlignment = -4108
$tmpl.VerticalAlignment = -4108
$tmpl.Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$tmpl.Clear()
$excel.CutCopyMode = $false

# --- Column widths -------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13.8
$ws.Columns.Item(2).ColumnWidth = 31.1
$ws.Columns.Item(3).ColumnWidth = 15.5
$ws.Columns.Item(4).ColumnWidth = 17.2

# --- Selection -------------------------------------------------------------
$ws.Range("E9").Select() | Out-Null
